{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line\n// that follows it, and the blank paragraph that separates them from the\n// bibliography text above. This mirrors the upstream Jekyll rebuild that\n// dropped the scraped site-chrome text from the generated document while\n// leaving the rest of the content (including the final blank paragraph and\n// the page-break paragraph) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex !== -1) {\n  const toDelete = [items[markerIndex]];\n  // The blank paragraph immediately preceding the marker paragraph is part\n  // of the block being removed.\n  if (markerIndex - 1 >= 0 && items[markerIndex - 1].text === \"\") {\n    toDelete.push(items[markerIndex - 1]);\n  }\n  // The paragraph right after the marker should be the copyright line.\n  if (\n    markerIndex + 1 < items.length &&\n    items[markerIndex + 1].text === targetTexts[1]\n  ) {\n    toDelete.push(items[markerIndex + 1]);\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the copyright/footer line\n# that follows it, and the blank paragraph that separates them from the\n# bibliography text above. This mirrors the upstream Jekyll rebuild that\n# dropped the scraped site-chrome text from the generated document while\n# leaving the rest of the content (including the final blank paragraph and\n# the page-break paragraph) untouched.\n\n$d = $word.ActiveDocument\n\n$markerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"$([char]0xA9) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$markerIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $markerText) {\n    $markerIndex = $i\n    break\n  }\n}\n\nif ($markerIndex -ge 1) {\n  $indexesToDelete = New-Object System.Collections.Generic.List[int]\n  [void]$indexesToDelete.Add($markerIndex)\n\n  if ($markerIndex - 1 -ge 1) {\n    $prevText = $d.Paragraphs.Item($markerIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($prevText -eq \"\") {\n      [void]$indexesToDelete.Add($markerIndex - 1)\n    }\n  }\n\n  if ($markerIndex + 1 -le $d.Paragraphs.Count) {\n    $nextText = $d.Paragraphs.Item($markerIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($nextText -eq $copyrightText) {\n      [void]$indexesToDelete.Add($markerIndex + 1)\n    }\n  }\n\n  $sorted = $indexesToDelete | Sort-Object -Descending\n  foreach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n  }\n}\n"}
